$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1 (we pass explicit values through Execute signature)
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

# --- Change 1 ------------------------------------------------------------
# " Sistema identificar Código del Usuario al este ejecutar una acción
# transaccional, mediante el inicio de sesión."
#   -> " Sistema buscar a un usuario mediante sus credenciales para poder
#       acceder al sistema."
$d.Content.Find.Execute(
    " Sistema identificar Código del Usuario al este ejecutar una acción transaccional, mediante el inicio de sesión.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Sistema buscar a un usuario mediante sus credenciales para poder acceder al sistema.",
    2) | Out-Null

# --- Change 2 --------------------------------------------------------------
# "Alumno" (Actores section) -> "Usuario"
$d.Content.Find.Execute(
    "Alumno",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Usuario",
    2) | Out-Null

# --- Change 3 ----------------------------------------------------------
# "intenta " + "loguearse" + " en el sistema." ->
# "intenta loguearse en el sistema seleccionando el botón “Ingresar” de la
#  interfaz Login."
$d.Content.Find.Execute(
    "intenta loguearse en el sistema.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "intenta loguearse en el sistema seleccionando el botón “Ingresar” de la interfaz Login.",
    2) | Out-Null

# --- Change 4 ----------------------------------------------------------
# "Busca alumno mediante " -> "Busca usuario mediante "
$d.Content.Find.Execute(
    "Busca alumno mediante ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Busca usuario mediante ",
    2) | Out-Null
